$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.407.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = "'1.713.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'224.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = "'0.5323"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").Value = "'1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").Value = "'0.2664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.78%  '
$ws.Range("D9").Value = "'0.06606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("D10").Value = "'20.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.02%  '
$ws.Range("D11").Value = "'0.07653"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = "'4.580"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.58%  '
$ws.Range("D13").Value = "'1.726.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").Value = "'1.946.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = "'0.5735"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.02%  '
$ws.Range("D16").Value = "'0.0₅8187"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = "'67.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("D18").Value = "'27.408.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.15%  '
$ws.Range("D19").Value = "'216.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.25%  '
$ws.Range("D20").Value = "'1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = "'4.680"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").Value = "'10.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.21%  '
$ws.Range("D23").Value = "'5.981"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.12%  '
$ws.Range("D24").Value = "'1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = "'1.771"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.54%  '
$ws.Range("D26").Value = "'141.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").Value = "'0.1219"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("D28").Value = "'7.289"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").Value = "'16.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.17%  '
$ws.Range("D30").Value = "'0.05426"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.39%  '
$ws.Range("D31").Value = "'1.297"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").Value = "'3.512"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.23%  '
$ws.Range("D33").Value = "'3.435"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("D34").Value = "'1.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("D35").Value = "'2.886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").Value = "'0.9509"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.46%  '
$ws.Range("D37").Value = "'2.423"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.12%  '
$ws.Range("D38").Value = "'0.5876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("D39").Value = "'0.01632"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.34%  '
$ws.Range("D40").Value = "'5.880"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = "'1.048.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").Value = "'1.007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").Value = "'0.8458"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = "'101.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").Value = "'1.853.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = "'58.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.24%  '
$ws.Range("D48").Value = "'0.4517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = "'1.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D50").Value = "'8.081"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("D51").Value = "'0.05250"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.39%  '
